$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Rename header labels: SMD -> effect, SE -> se
$ws.Range("H1").Value = "effect"

# Add the new "effect=smd" annotation in K10 (before I1's "se" so the new
# shared-string entries land in the order effect, effect=smd, se)
$ws.Range("K10").Value = "effect=smd"

$ws.Range("I1").Value = "se"

# Update the selected cell to H1
$ws.Range("H1").Select()
